# Commit: "Tests: updated unit tests after adding linprog support"
#
# The "general" settings sheet gains a new "LP solver (linprog or gurobi)"
# row (value "gurobi") right after the existing "NLP solver" row. Every
# row that used to follow shifts down by one; all other sheets are
# untouched (their apparent diffs are just shared-string / style index
# renumbering caused by this single insertion).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

# Push "Number of exp. conditions" (and everything below it) down one row.
$ws.Rows.Item(5).Insert()

# Fill in the new row.
$ws.Cells.Item(5, 1).Value = "LP solver (linprog or gurobi)"
$ws.Cells.Item(5, 2).Value = "gurobi"

# The label cells in column A use a bold/bordered style; copy that look
# from the row above ("NLP solver") onto the new label cell, then switch
# its horizontal alignment to left (the other label cells stay centered).
$ws.Cells.Item(4, 1).Copy()
$ws.Cells.Item(5, 1).PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Cells.Item(5, 1).HorizontalAlignment = -4131  # xlLeft

# Match the row height used elsewhere in this sheet for this style.
$ws.Rows.Item(5).RowHeight = 13.8
